$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.378.65'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.55'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.59'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6294'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07611'
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.62'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07752'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6816'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001050'
$ws.Range('E14').Value = '  -4.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.16'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.125'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.386.64'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.02'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.468'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '159.18'
$ws.Range('E23').Value = '  +1.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.1392'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.447'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.71'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.448'
$ws.Range('E27').Value = '  +9.60%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05640'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.112'
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.058'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.833'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.158'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7015'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.586'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01825'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.238.09'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.728'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.421'
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9063'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.49'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.63'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.156'
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000118'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4002'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.052'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1157'
$ws.Range('E48').Value = '  +2.08%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.685'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05703'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4628'
$ws.Range('E51').Value = '  -0.12%  '
